# Reorders the data rows (2-23) of Sheet1 into the new "kobo adjustment" order.
# The set of (file, r_count, currency, sum) tuples is unchanged; only the row
# order (and therefore which row each file/currency/sum combination sits in)
# is updated. Column E ("built_in_total") holds the sum formatted as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row order / values for A2:E23
# file, r_count, currency, sum (numeric), built_in_total (text)
$rows = @(
    @("87811004_0322_RO", 423,  "RON", 9983.42,    "9983.42"),
    @("87811004_0322_PE", 15,   "PEN", 141.96,     "141.96"),
    @("87811004_0322_HU", 954,  "HUF", 1882603,    "1882603"),
    @("87811004_0322_EU", 716,  "EUR", 3090.3,     "3090.3"),
    @("87811004_0322_MX", 75,   "MXN", 4825.1,     "4825.1"),
    @("87811004_0322_LL", 37,   "USD", 96.6,       "96.6"),
    @("87811004_0322_BG", 6,    "BGN", 22.44,      "22.44"),
    @("87811004_0322_BR", 37,   "BRL", 334.81,     "334.81"),
    @("87811004_0322_CA", 324,  "CAD", 1753.5,     "1753.5"),
    @("87811004_0322_CZ", 12,   "CZK", 636.36,     "636.36"),
    @("87811004_0322_CL", 31,   "CLP", 63087,      "63087"),
    @("87811004_0322_CO", 31,   "COP", 286230,     "286230"),
    @("87811004_0322_NZ", 51,   "NZD", 208.33,     "208.33"),
    @("87811004_0322_AU", 367,  "AUD", 3178.64,    "3178.64"),
    @("87811004_0322_CH", 79,   "CHF", 243.39,     "243.39"),
    @("87811004_0322_NO", 26,   "NOK", 682.5,      "682.5"),
    @("87811004_0322_US", 1487, "USD", 10299.1,    "10299.1"),
    @("87811004_0322_DK", 21,   "DKK", 268.24,     "268.24"),
    @("87811004_0322_PL", 46,   "PLN", 517.94,     "517.94"),
    @("87811004_0322_SE", 32,   "SEK", 931.77,     "931.77"),
    @("87811004_0322_JP", 23,   "JPY", 5950,       "5950"),
    @("87811004_0322_GB", 477,  "GBP", 1593.33,    "1593.33")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Force column E to be stored as text (matching the original file's
    # shared-string "built_in_total" values) rather than being auto-coerced
    # into a number, then drop the quote-prefix style that Excel applies so
    # the cell keeps its original (default) style.
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = "'" + $row[4]
    $eCell.Style = "Normal"

    $r++
}
